$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 840.75
$ws.Cells.Item(12, 9).Value = 787.5
$ws.Cells.Item(12, 11).Value = 787.5
$ws.Cells.Item(12, 13).Value = -617.5
$ws.Cells.Item(18, 8).Value = 980.8182
$ws.Cells.Item(18, 9).Value = 974.125
$ws.Cells.Item(18, 10).Value = 998.6667
$ws.Cells.Item(18, 11).Value = 974.125
$ws.Cells.Item(18, 12).Value = 998.6667
$ws.Cells.Item(18, 13).Value = -690.125
$ws.Cells.Item(18, 14).Value = -1566.6667
$ws.Cells.Item(32, 8).Value = 794.6667
$ws.Cells.Item(32, 10).Value = 794.6667
$ws.Cells.Item(32, 12).Value = 794.6667
$ws.Cells.Item(32, 14).Value = -1446.6667
$ws.Cells.Item(40, 8).Value = 7062.7144
$ws.Cells.Item(40, 10).Value = 8332.666999999999
$ws.Cells.Item(40, 12).Value = 8332.666999999999
$ws.Cells.Item(40, 14).Value = -8682.666999999999
$ws.Cells.Item(41, 8).Value = 827.4167
$ws.Cells.Item(41, 9).Value = 629.9091
$ws.Cells.Item(41, 11).Value = 629.9091
$ws.Cells.Item(41, 13).Value = -189.9091
$ws.Cells.Item(80, 8).Value = 335.8
$ws.Cells.Item(80, 10).Value = 344
$ws.Cells.Item(80, 12).Value = 1032
$ws.Cells.Item(80, 14).Value = -3028
$ws.Cells.Item(83, 8).Value = 335.8
$ws.Cells.Item(83, 10).Value = 344
$ws.Cells.Item(83, 12).Value = 3096
$ws.Cells.Item(83, 14).Value = -13080
$ws.Cells.Item(87, 8).Value = 92451.664
$ws.Cells.Item(87, 10).Value = 92451.664
$ws.Cells.Item(87, 12).Value = 92451.664
$ws.Cells.Item(87, 14).Value = -94947.664
$ws.Cells.Item(90, 8).Value = 92451.664
$ws.Cells.Item(90, 10).Value = 92451.664
$ws.Cells.Item(90, 12).Value = 277354.992
$ws.Cells.Item(90, 14).Value = -289834.992

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(56, 8).Value = 0
$ws.Cells.Item(56, 9).Value = 0
$ws.Cells.Item(56, 11).Value = 0
$ws.Cells.Item(56, 13).ClearContents()
$ws.Cells.Item(61, 8).Value = 2805.625
$ws.Cells.Item(61, 9).Value = 2267.6667
$ws.Cells.Item(61, 10).Value = 4419.5
$ws.Cells.Item(61, 11).Value = 2267.6667
$ws.Cells.Item(61, 12).Value = 4419.5
$ws.Cells.Item(61, 13).Value = -2055.6667
$ws.Cells.Item(61, 14).Value = -4843.5
$ws.Cells.Item(136, 8).Value = 2805.625
$ws.Cells.Item(136, 9).Value = 2267.6667
$ws.Cells.Item(136, 10).Value = 4419.5
$ws.Cells.Item(136, 11).Value = 6803.000100000001
$ws.Cells.Item(136, 12).Value = 13258.5
$ws.Cells.Item(136, 13).Value = -4253.000100000001
$ws.Cells.Item(136, 14).Value = -18358.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1414.5
$ws.Cells.Item(20, 9).Value = 697.6
$ws.Cells.Item(20, 11).Value = 697.6
$ws.Cells.Item(20, 13).Value = -450.6
$ws.Cells.Item(22, 8).Value = 345
$ws.Cells.Item(22, 9).Value = 345
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 345
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -172
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(86, 8).Value = 3513.1875
$ws.Cells.Item(86, 10).Value = 4649.9
$ws.Cells.Item(86, 12).Value = 4649.9
$ws.Cells.Item(86, 14).Value = -6895.9
$ws.Cells.Item(89, 8).Value = 3513.1875
$ws.Cells.Item(89, 10).Value = 4649.9
$ws.Cells.Item(89, 12).Value = 23249.5
$ws.Cells.Item(89, 14).Value = -34481.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4394.5557
$ws.Cells.Item(31, 9).Value = 1447.0625
$ws.Cells.Item(31, 11).Value = 1447.0625
$ws.Cells.Item(31, 13).Value = -1152.0625
$ws.Cells.Item(34, 8).Value = 4394.5557
$ws.Cells.Item(34, 9).Value = 1447.0625
$ws.Cells.Item(34, 11).Value = 1447.0625
$ws.Cells.Item(34, 13).Value = -1245.0625
$ws.Cells.Item(62, 8).Value = 3249.5
$ws.Cells.Item(62, 10).Value = 3249.5
$ws.Cells.Item(62, 12).Value = 3249.5
$ws.Cells.Item(62, 14).Value = -4497.5
$ws.Cells.Item(65, 8).Value = 3249.5
$ws.Cells.Item(65, 10).Value = 3249.5
$ws.Cells.Item(65, 12).Value = 16247.5
$ws.Cells.Item(65, 14).Value = -22487.5
$ws.Cells.Item(87, 8).Value = 30000
$ws.Cells.Item(87, 10).Value = 30000
$ws.Cells.Item(87, 12).Value = 30000
$ws.Cells.Item(87, 14).Value = -32372
$ws.Cells.Item(90, 8).Value = 30000
$ws.Cells.Item(90, 10).Value = 30000
$ws.Cells.Item(90, 12).Value = 90000
$ws.Cells.Item(90, 14).Value = -101856
$ws.Cells.Item(99, 8).Value = 2874.75
$ws.Cells.Item(99, 9).Value = 2874.75
$ws.Cells.Item(99, 11).Value = 2874.75
$ws.Cells.Item(99, 13).Value = -1376.75
$ws.Cells.Item(126, 8).Value = 2874.75
$ws.Cells.Item(126, 9).Value = 2874.75
$ws.Cells.Item(126, 11).Value = 8624.25
$ws.Cells.Item(126, 13).Value = -6154.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 133.76923
$ws.Cells.Item(12, 10).Value = 196.625
$ws.Cells.Item(12, 12).Value = 589.875
$ws.Cells.Item(12, 14).Value = -935.875
$ws.Cells.Item(60, 8).Value = 1249.6875
$ws.Cells.Item(60, 10).Value = 2268.75
$ws.Cells.Item(60, 12).Value = 6806.25
$ws.Cells.Item(60, 14).Value = -7308.25
$ws.Cells.Item(122, 8).Value = 668.8
$ws.Cells.Item(122, 9).Value = 372.75
$ws.Cells.Item(122, 10).Value = 866.1667
$ws.Cells.Item(122, 11).Value = 3354.75
$ws.Cells.Item(122, 12).Value = 7795.5003
$ws.Cells.Item(122, 13).Value = -904.75
$ws.Cells.Item(122, 14).Value = -12695.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 343.5
$ws.Cells.Item(2, 9).Value = 153.27272
$ws.Cells.Item(2, 11).Value = 153.27272
$ws.Cells.Item(2, 13).Value = -40.27271999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1244.6364
$ws.Cells.Item(16, 9).Value = 1244.6364
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 1244.6364
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = -1074.6364
$ws.Cells.Item(16, 14).ClearContents()
$ws.Cells.Item(22, 8).Value = 711.2308
$ws.Cells.Item(22, 9).Value = 182.66667
$ws.Cells.Item(22, 10).Value = 869.8
$ws.Cells.Item(22, 11).Value = 182.66667
$ws.Cells.Item(22, 12).Value = 869.8
$ws.Cells.Item(22, 13).Value = 112.33333
$ws.Cells.Item(22, 14).Value = -1459.8
$ws.Cells.Item(27, 8).Value = 711.2308
$ws.Cells.Item(27, 9).Value = 182.66667
$ws.Cells.Item(27, 10).Value = 869.8
$ws.Cells.Item(27, 11).Value = 182.66667
$ws.Cells.Item(27, 12).Value = 869.8
$ws.Cells.Item(27, 13).Value = -75.66667000000001
$ws.Cells.Item(27, 14).Value = -1083.8
$ws.Cells.Item(46, 8).Value = 6254.95
$ws.Cells.Item(46, 9).Value = 5287.375
$ws.Cells.Item(46, 10).Value = 6900
$ws.Cells.Item(46, 11).Value = 5287.375
$ws.Cells.Item(46, 12).Value = 6900
$ws.Cells.Item(46, 13).Value = -5099.375
$ws.Cells.Item(46, 14).Value = -7276
$ws.Cells.Item(61, 8).Value = 90914160
$ws.Cells.Item(61, 9).Value = 142861330
$ws.Cells.Item(61, 10).Value = 6624.75
$ws.Cells.Item(61, 11).Value = 142861330
$ws.Cells.Item(61, 12).Value = 6624.75
$ws.Cells.Item(61, 13).Value = -142861128
$ws.Cells.Item(61, 14).Value = -7028.75
$ws.Cells.Item(93, 8).Value = 1836.7858
$ws.Cells.Item(93, 9).Value = 1939.375
$ws.Cells.Item(93, 10).Value = 1700
$ws.Cells.Item(93, 11).Value = 1939.375
$ws.Cells.Item(93, 12).Value = 1700
$ws.Cells.Item(93, 13).Value = -691.375
$ws.Cells.Item(93, 14).Value = -4196
$ws.Cells.Item(113, 8).Value = 90914160
$ws.Cells.Item(113, 9).Value = 142861330
$ws.Cells.Item(113, 10).Value = 6624.75
$ws.Cells.Item(113, 11).Value = 142861330
$ws.Cells.Item(113, 12).Value = 6624.75
$ws.Cells.Item(113, 13).Value = -142859160
$ws.Cells.Item(113, 14).Value = -10964.75
$ws.Cells.Item(132, 8).Value = 3002
$ws.Cells.Item(132, 9).Value = 3002
$ws.Cells.Item(132, 11).Value = 9006
$ws.Cells.Item(132, 13).Value = -6476

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 37037816
$ws.Cells.Item(107, 9).Value = 41667388
$ws.Cells.Item(107, 11).Value = 125002164
$ws.Cells.Item(107, 13).Value = -125000244
$ws.Cells.Item(122, 8).Value = 2838.75
$ws.Cells.Item(122, 9).Value = 2501.6667
$ws.Cells.Item(122, 11).Value = 7505.000100000001
$ws.Cells.Item(122, 13).Value = -5055.000100000001
$ws.Cells.Item(126, 8).Value = 6682.8
$ws.Cells.Item(126, 10).Value = 6682.8
$ws.Cells.Item(126, 12).Value = 20048.4
$ws.Cells.Item(126, 14).Value = -24988.4
